# Auto-generated edit script: applies the row reordering + odds updates
# described by the commit diff for "Chile Primera Division.xlsx".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

function Swap-Rows($rowA, $rowB) {
    foreach ($c in $cols) {
        $refA = "$c$rowA"
        $refB = "$c$rowB"
        $valA = $ws.Range($refA).Value2
        $valB = $ws.Range($refB).Value2
        $ws.Range($refA).Value2 = $valB
        $ws.Range($refB).Value2 = $valA
    }
}

function Cycle-Rows($rowList) {
    # Each row in the list takes on the values that were originally in the
    # PREVIOUS row of the list (wrapping around), i.e. for (r0, r1, r2, ...):
    #   new(r0) = old(r_last), new(r1) = old(r0), new(r2) = old(r1), ...
    $snapshot = @{}
    foreach ($r in $rowList) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Range("$c$r").Value2
        }
        $snapshot[$r] = $rowVals
    }
    $n = $rowList.Length
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $rowList[$i]
        $srcRow = $rowList[(($i - 1) + $n) % $n]
        foreach ($c in $cols) {
            $ws.Range("$c$destRow").Value2 = $snapshot[$srcRow][$c]
        }
    }
}

# --- Row pairs that simply swapped places (match ordering changed) ---
Swap-Rows 64 65
Swap-Rows 112 114
Swap-Rows 115 117
Swap-Rows 121 122

# --- Three-way rotation: new(118) = old(120), new(119) = old(118), new(120) = old(119) ---
Cycle-Rows @(118, 119, 120)

# --- Standalone odds recalculations (no row movement) ---
# Row 201
$ws.Range("M201").Value2 = 1.571
$ws.Range("N201").Value2 = 4
$ws.Range("O201").Value2 = 6
$ws.Range("Q201").Value2 = 2
$ws.Range("R201").Value2 = 1.85
$ws.Range("T201").Value2 = 2.05
$ws.Range("U201").Value2 = 1.8

# Row 202
$ws.Range("M202").Value2 = 2.2
$ws.Range("N202").Value2 = 3.25
$ws.Range("O202").Value2 = 3.25
$ws.Range("Q202").Value2 = 1.9
$ws.Range("R202").Value2 = 1.95
$ws.Range("T202").Value2 = 2.025
$ws.Range("U202").Value2 = 1.825

# Row 203
$ws.Range("S203").Value2 = 2.5
$ws.Range("T203").Value2 = 1.825
$ws.Range("U203").Value2 = 2.025

# Row 204
$ws.Range("S204").Value2 = 2.75
$ws.Range("T204").Value2 = 2.05
$ws.Range("U204").Value2 = 1.8

# Row 205
$ws.Range("M205").Value2 = 2.55
$ws.Range("O205").Value2 = 2.875
$ws.Range("P205").Value2 = 0
$ws.Range("Q205").Value2 = 1.8
$ws.Range("R205").Value2 = 2.05

# Row 206
$ws.Range("T206").Value2 = 1.825
$ws.Range("U206").Value2 = 2.025

# Row 207
$ws.Range("M207").Value2 = 2.875
$ws.Range("O207").Value2 = 2.375
$ws.Range("Q207").Value2 = 1.775
$ws.Range("R207").Value2 = 2.1
$ws.Range("T207").Value2 = 1.975
$ws.Range("U207").Value2 = 1.875
